# The experiment's trial table is restructured: the single stimulus/probe
# row pair (with a "change" flag column) becomes two side-by-side blocks -
# one for the initial stimulus array (stim1_x/stim1_y/sitm1_color) and one
# for the probe (probe1_x/probe1_y/probe1_color) - plus an "answer" column
# and a "probe1" fixation-style label, all collapsed onto a single data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: new headers
$ws.Range("A1").Value = "stim1_x"
$ws.Range("B1").Value = "stim1_y"
$ws.Range("C1").Value = "sitm1_color"
$ws.Range("D1").Value = "probe1_x"
$ws.Range("E1").Value = "probe1_y"
$ws.Range("F1").Value = "probe1_color"
$ws.Range("G1").Value = "answer"

# Row 2: the merged stimulus + probe trial data
$ws.Range("A2").Value = -0.25
$ws.Range("B2").Value = -0.25
$ws.Range("C2").Value = "black"
$ws.Range("D2").Value = -0.25
$ws.Range("E2").Value = -0.25
$ws.Range("F2").Value = "white"
$ws.Range("G2").Value = "probe1"

# The old third row is no longer needed now that everything fits on row 2
$ws.Range("A3:D3").ClearContents() | Out-Null

# Column F ("probe1_color") is now the widest header/value, so it gets
# auto-fit to its content width
$ws.Columns.Item(6).ColumnWidth = 11

# Leave the selection where the user's cursor ended up after the edit
$ws.Range("B11").Select() | Out-Null
